# Update tokens file and get tokens from the token file
#
# Target shape (per diff):
#   Row1: userEmail | longToken | tokenExpiry               (unchanged)
#   Row2: user_123  | <long fb token #1>        | 1730486029253  (ms-epoch int, General fmt, wrap)
#   Row3: (blank)   | <long fb token #2>        | 1730491898045  (ms-epoch int, General fmt, wrap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$token1 = "EAB3k2pU1ZBxwBO6ixHJeIO6kJu3wGOtPkyrUxA8zRL4ZAZAPrS6oLq1Nn0eXyDRKpc9aNO7cr9KsG7ZCAs23aTj7tlnwaZC3Mu5gOwZBRMQ2vRZCQK48vjx3L35Pm7biuICrf1s39CZAmExTlZBdOjJE3fxbtcMHwKIZBbDehr3hqdaeMvLp46jklc3UyB"
$token2 = "EAAYEnKKz8boBO5oswUfODAELXCONlHCM3UDCBTcvWZBSUZAelyM17CUVW04DPHM7ZAya1NEgSipigLu32zHKXYVBuUVgGGay7SBzniPuvSYxzRFBV6VtmVzecB73CulCqK2ptnMGm7zJRifdRoQCbc3dg8c17AzoLgOwX8ZAYJWAsf7SdvB2lbndDUJZCnzyZAYkrjVZA6R"

# --- Row 2: replace the old short token string and refresh the timestamp ---
$ws.Range("B2").Value = $token1
$ws.Range("C2").Value = 1730486029253

# --- Row 3: brand new row holding a second token + its own expiry timestamp ---
$ws.Range("B3").Value = $token2
$ws.Range("C3").Value = 1730491898045

# --- Formatting: the timestamp column no longer uses a date format, and the
#     token / timestamp cells should wrap their (very long) text ---
$fmtRange = $ws.Range("B2:C3")
$fmtRange.NumberFormat = "General"
$fmtRange.WrapText = $true

# --- Row heights so the wrapped, multi-line text is visible ---
$ws.Rows.Item(2).RowHeight = 79.85
$ws.Rows.Item(3).RowHeight = 58.2

# --- Column widths so the long tokens have room ---
# (the engine quantizes ColumnWidth to ~1/6-character steps, so these values
#  are chosen to land on the raw stored width closest to the 59.1 / 60.51
#  targets: 58.3 -> ~59.17, 59.6 -> ~60.5)
$ws.Columns.Item(2).ColumnWidth = 58.3
$ws.Columns.Item(3).ColumnWidth = 59.6

# --- Active cell / selection moves to the newly added C3 ---
$null = $ws.Range("C3").Select()
